# "Error Calculations and Plots"
# Reconciles the missing_data.xlsx sample against the reference measurements:
#  - fills back in a handful of values that were previously blanked out
#  - blanks out a handful of values that were previously present
#  - drops two rows (RM 232 and SC 92) that are no longer part of the sample
#  - fills in the previously-missing "B" value for SC 193

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Single-cell value corrections -----------------------------------------

# RM 8: column B value restored
$ws.Range("C3").Value = 11.2

# RM 9: column D value now missing
$ws.Range("E4").ClearContents()

# RM 14: column B value now missing
$ws.Range("C5").ClearContents()

# RM 42: column D value restored
$ws.Range("E9").Value = -6.8

# RM 52 a: column D value restored
$ws.Range("E10").Value = -6.1

# RM 88: column D value now missing
$ws.Range("E13").ClearContents()

# RM 90: column D value now missing
$ws.Range("E14").ClearContents()

# RM 135: column B value restored
$ws.Range("C21").Value = 12.7

# RM 140: column B value now missing
$ws.Range("C23").ClearContents()

# --- Row removals ------------------------------------------------------------
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(28).EntireRow.Delete()   # SC 92
$ws.Rows.Item(26).EntireRow.Delete()   # RM 232

# --- Fill in the previously-missing value for SC 193 (now row 32) -----------
$ws.Range("C32").Value = 10.5
